$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.555.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.850.69"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9973"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6306"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9990"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07498"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2918"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.71"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07740"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.850.73"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.025"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6834"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001046"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.271"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.546.25"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "230.00"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9987"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.581"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9991"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.46"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.539"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.28%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.57"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06534"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +15.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.433"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.485"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.10%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.107"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.844"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.24%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6999"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.576"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01863"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.266.63"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.841"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.862"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9375"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.64%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "RocketPoolETH"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.032.51"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.30%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9996"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.33"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.741"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.120"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1168"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.56%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.019"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.51%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3969"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.69%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000115"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.55%  "
